# Atualizacao rapida de agenda as 15:27:05,40
#
# Inserts a new agenda row (new order at "Cetep") above the current row 8,
# shifting every row from 8..158 down to 9..159, updates row 7 with a new
# observation + status, and fixes up the selection / row heights.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift rows 8..158 down to 9..159 (preserves values + styles) -------
$ws.Range("A8:H158").Copy() | Out-Null
$ws.Range("A9:H159").PasteSpecial() | Out-Null

# The paste above can leave the brand-new trailing row (159) without its
# cells materialized (they are visually blank). Force them into existence
# and give them row 158's formatting (blank cells, same style pattern as
# every other untouched agenda row) so the sheet dimension grows correctly.
for ($c = 1; $c -le 8; $c++) {
    $ws.Cells.Item(159, $c).Value = ""
}
$ws.Range("A158:H158").Copy() | Out-Null
$ws.Range("A159:H159").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 2. Update row 7: new observation + status -----------------------------
$ws.Range("E7").Value = "Giovani esteve no local mas o pessoal estava usando as máquinas e ninguém pode ir lá abrir pra ele o local onde fica a central."
$ws.Range("G7").Value = "Concluido"
$ws.Rows.Item(7).RowHeight = 30

# --- 3. Fill the newly inserted row 8 with the new agenda entry ------------
$ws.Range("A8").Value = "Giovani"
$ws.Range("B8").Value = "'0790"
$ws.Range("C8").Value = "Cetep"
$ws.Range("D8").Value = "Solicitado pelo cliente, passar pro DDNS"
$ws.Range("E8").Value = "Essa ordem estava pra amanhã.. Giovani aparentemente esteve no local e combinou com a diretora a visita pra amanhã, visto que ele solicitou um técnico da claro também no colégio pra fazer a execução do serviço."
$ws.Range("G8").Value = "Concluido"
$ws.Rows.Item(8).RowHeight = 45

# --- 4. Restore the active cell selection -----------------------------------
$ws.Range("H11").Select() | Out-Null
